$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Btc"
$ws.Range("C2").Value = "Erbb2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1369323333333334
$ws.Range("H2").Value = 0.410797
$ws.Range("I2").Value = 0.08243498013860336
$ws.Range("J2").Value = 0.08243498013860337
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.885873333333334
$ws.Range("N2").Value = 8.657620000000001
$ws.Range("O2").Value = 0.3070415651026022
$ws.Range("P2").Value = 0.3070415651026022
$ws.Range("Q2").Value = 0.3951693692377779
$ws.Range("R2").Value = 3.556524323140001
$ws.Range("S2").Value = 0.0253109653209587
$ws.Range("T2").Value = 0.02531096532095871

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Btc"
$ws.Range("C3").Value = "Erbb2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1369323333333334
$ws.Range("H3").Value = 0.410797
$ws.Range("I3").Value = 0.08243498013860336
$ws.Range("J3").Value = 0.08243498013860337
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.165953666666667
$ws.Range("N3").Value = 9.497861
$ws.Range("O3").Value = 0.3368406220840099
$ws.Range("P3").Value = 0.3368406220840099
$ws.Range("Q3").Value = 0.4335214228018889
$ws.Range("R3").Value = 3.901692805217
$ws.Range("S3").Value = 0.02776744999137016
$ws.Range("T3").Value = 0.02776744999137016

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Btc"
$ws.Range("C4").Value = "Erbb2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1369323333333334
$ws.Range("H4").Value = 0.410797
$ws.Range("I4").Value = 0.08243498013860336
$ws.Range("J4").Value = 0.08243498013860337
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.327024333333334
$ws.Range("N4").Value = 9.981073
$ws.Range("O4").Value = 0.3539776838580724
$ws.Range("P4").Value = 0.3539776838580724
$ws.Range("Q4").Value = 0.4555772050201112
$ws.Range("R4").Value = 4.100194845181
$ws.Range("S4").Value = 0.02918014333834902
$ws.Range("T4").Value = 0.02918014333834902

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Btc"
$ws.Range("C5").Value = "Erbb2"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.1369323333333334
$ws.Range("H5").Value = 0.410797
$ws.Range("I5").Value = 0.08243498013860336
$ws.Range("J5").Value = 0.08243498013860337
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.020115
$ws.Range("N5").Value = 0.060345
$ws.Range("O5").Value = 0.002140128955315263
$ws.Range("P5").Value = 0.002140128955315263
$ws.Range("Q5").Value = 0.002754393885
$ws.Range("R5").Value = 0.024789544965
$ws.Range("S5").Value = 0.0001764214879254637
$ws.Range("T5").Value = 0.0001764214879254637

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Btc"
$ws.Range("C6").Value = "Erbb2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.524162666666667
$ws.Range("H6").Value = 4.572488
$ws.Range("I6").Value = 0.9175650198613966
$ws.Range("J6").Value = 0.9175650198613967
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.885873333333334
$ws.Range("N6").Value = 8.657620000000001
$ws.Range("O6").Value = 0.3070415651026022
$ws.Range("P6").Value = 0.3070415651026022
$ws.Range("Q6").Value = 4.398540395395556
$ws.Range("R6").Value = 39.58686355856
$ws.Range("S6").Value = 0.2817305997816435
$ws.Range("T6").Value = 0.2817305997816436

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Btc"
$ws.Range("C7").Value = "Erbb2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.524162666666667
$ws.Range("H7").Value = 4.572488
$ws.Range("I7").Value = 0.9175650198613966
$ws.Range("J7").Value = 0.9175650198613967
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.165953666666667
$ws.Range("N7").Value = 9.497861
$ws.Range("O7").Value = 0.3368406220840099
$ws.Range("P7").Value = 0.3368406220840099
$ws.Range("Q7").Value = 4.825428383129777
$ws.Range("R7").Value = 43.428855448168
$ws.Range("S7").Value = 0.3090731720926398
$ws.Range("T7").Value = 0.3090731720926398

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Btc"
$ws.Range("C8").Value = "Erbb2"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.524162666666667
$ws.Range("H8").Value = 4.572488
$ws.Range("I8").Value = 0.9175650198613966
$ws.Range("J8").Value = 0.9175650198613967
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.327024333333334
$ws.Range("N8").Value = 9.981073
$ws.Range("O8").Value = 0.3539776838580724
$ws.Range("P8").Value = 0.3539776838580724
$ws.Range("Q8").Value = 5.070926279958222
$ws.Range("R8").Value = 45.638336519624
$ws.Range("S8").Value = 0.3247975405197234
$ws.Range("T8").Value = 0.3247975405197234

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Btc"
$ws.Range("C9").Value = "Erbb2"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.524162666666667
$ws.Range("H9").Value = 4.572488
$ws.Range("I9").Value = 0.9175650198613966
$ws.Range("J9").Value = 0.9175650198613967
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.020115
$ws.Range("N9").Value = 0.060345
$ws.Range("O9").Value = 0.002140128955315263
$ws.Range("P9").Value = 0.002140128955315263
$ws.Range("Q9").Value = 0.03065853204
$ws.Range("R9").Value = 0.27592678836
$ws.Range("S9").Value = 0.001963707467389799
$ws.Range("T9").Value = 0.0019637074673898
